$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Match-ID renumbering: rows 45-47 (was M001) become M003,
#    and rows 82-84 (was M003) become M001 - they swap identifiers.
# ---------------------------------------------------------------------------
$ws.Range("A45").Value = "M003"
$ws.Range("A46").Value = "M003"
$ws.Range("A47").Value = "M003"

$ws.Range("A82").Value = "M001"
$ws.Range("A83").Value = "M001"
$ws.Range("A84").Value = "M001"

# ---------------------------------------------------------------------------
# 2) Rows 70-75 become two new manually matched 3-row blocks (M028 / M029),
#    formatted exactly like existing "Dr / Receipt" 3-row blocks elsewhere in
#    the ledger (rows 82-84 and 79-81 use the same Dr/Receipt layout with the
#    amount in column K). Copy formatting from those reference blocks first,
#    then set the new values/text.
# ---------------------------------------------------------------------------
$ws.Range("A82:L84").Copy() | Out-Null
$ws.Range("A70").PasteSpecial(-4122) | Out-Null

$ws.Range("A79:L81").Copy() | Out-Null
$ws.Range("A73").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Row 70 (M028 header line) ---
$ws.Range("A70").Value = "M028"
$ws.Range("B70").Value = ""
$ws.Range("C70").Value = "24/Dec/2025"
$ws.Range("D70").Value = "Dr"
$ws.Range("E70").Value = "Brac Bank PLC-CD-A/C-2028701210002"
$ws.Range("F70").Value = ""
$ws.Range("G70").Value = ""
$ws.Range("H70").Value = "Receipt"
$ws.Range("I70").Value = "67"
$ws.Range("J70").Value = ""
$ws.Range("K70").Value = "45000000"
$ws.Range("L70").Value = "Manual"

# --- Row 71 (M028 match detail line) ---
$ws.Range("A71").Value = "M028"
$ws.Range("B71").Value = "Manual Match`nLender Amount: 45000000.00`nBorrower Amount: 45000000.00"
$ws.Range("C71").Value = ""
$ws.Range("D71").Value = ""
$ws.Range("E71").Value = "Interunit Funs Transfer as Interunit Loan A/C-Steel Unit, MTB# 1105"
$ws.Range("F71").Value = ""
$ws.Range("G71").Value = ""
$ws.Range("H71").Value = ""
$ws.Range("I71").Value = ""
$ws.Range("J71").Value = ""
$ws.Range("K71").Value = ""
$ws.Range("L71").Value = "Manual"

# --- Row 72 (M028 entered-by line) ---
$ws.Range("A72").Value = "M028"
$ws.Range("B72").Value = ""
$ws.Range("C72").Value = ""
$ws.Range("D72").Value = "Entered By :"
$ws.Range("E72").Value = "ashiq"
$ws.Range("F72").Value = ""
$ws.Range("G72").Value = ""
$ws.Range("H72").Value = ""
$ws.Range("I72").Value = ""
$ws.Range("J72").Value = ""
$ws.Range("K72").Value = ""
$ws.Range("L72").Value = "Manual"

# --- Row 73 (M029 header line) ---
$ws.Range("A73").Value = "M029"
$ws.Range("B73").Value = ""
$ws.Range("C73").Value = "28/Dec/2025"
$ws.Range("D73").Value = "Dr"
$ws.Range("E73").Value = "Brac Bank PLC-CD-A/C-2028701210002"
$ws.Range("F73").Value = ""
$ws.Range("G73").Value = ""
$ws.Range("H73").Value = "Receipt"
$ws.Range("I73").Value = "68"
$ws.Range("J73").Value = ""
$ws.Range("K73").Value = "100000"
$ws.Range("L73").Value = "Manual"

# --- Row 74 (M029 match detail line) ---
$ws.Range("A74").Value = "M029"
$ws.Range("B74").Value = "Manual Match`nLender Amount: 100000.00`nBorrower Amount: 100000.00"
$ws.Range("C74").Value = ""
$ws.Range("D74").Value = ""
$ws.Range("E74").Value = "Interunit Funs Transfer as Interunit Loan A/C-Steel Unit, MTB# 1105"
$ws.Range("F74").Value = ""
$ws.Range("G74").Value = ""
$ws.Range("H74").Value = ""
$ws.Range("I74").Value = ""
$ws.Range("J74").Value = ""
$ws.Range("K74").Value = ""
$ws.Range("L74").Value = "Manual"

# --- Row 75 (M029 entered-by line) ---
$ws.Range("A75").Value = "M029"
$ws.Range("B75").Value = ""
$ws.Range("C75").Value = ""
$ws.Range("D75").Value = "Entered By :"
$ws.Range("E75").Value = "ashiq"
$ws.Range("F75").Value = ""
$ws.Range("G75").Value = ""
$ws.Range("H75").Value = ""
$ws.Range("I75").Value = ""
$ws.Range("J75").Value = ""
$ws.Range("K75").Value = ""
$ws.Range("L75").Value = "Manual"
